$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44245
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 800
$ws.Cells.Item(2, 11).Value = 850
$ws.Cells.Item(2, 12).Value = 900
$ws.Cells.Item(2, 13).Value = 875
$ws.Cells.Item(2, 16).Value = 875

$ws.Cells.Item(3, 4).Value = 44245
$ws.Cells.Item(3, 9).Value = "Segunda"
$ws.Cells.Item(3, 10).Value = 1000
$ws.Cells.Item(3, 11).Value = 750
$ws.Cells.Item(3, 12).Value = 800
$ws.Cells.Item(3, 13).Value = 775
$ws.Cells.Item(3, 16).Value = 775

$ws.Cells.Item(4, 4).Value = 44224
$ws.Cells.Item(4, 9).Value = "Segunda"
$ws.Cells.Item(4, 10).Value = 800
$ws.Cells.Item(4, 11).Value = 850
$ws.Cells.Item(4, 12).Value = 900
$ws.Cells.Item(4, 13).Value = 875
$ws.Cells.Item(4, 16).Value = 875

$ws.Cells.Item(5, 4).Value = 44935
$ws.Cells.Item(5, 9).Value = "Segunda"
$ws.Cells.Item(5, 10).Value = 1000
$ws.Cells.Item(5, 11).Value = 400
$ws.Cells.Item(5, 12).Value = 500
$ws.Cells.Item(5, 13).Value = 460
$ws.Cells.Item(5, 16).Value = 460

$ws.Cells.Item(6, 4).Value = 44474
$ws.Cells.Item(6, 9).Value = "Segunda"
$ws.Cells.Item(6, 10).Value = 200
$ws.Cells.Item(6, 11).Value = 600
$ws.Cells.Item(6, 12).Value = 700
$ws.Cells.Item(6, 13).Value = 650
$ws.Cells.Item(6, 16).Value = 650

$ws.Cells.Item(7, 4).Value = 44874
$ws.Cells.Item(7, 9).Value = "Tercera"
$ws.Cells.Item(7, 10).Value = 1200
$ws.Cells.Item(7, 11).Value = 450
$ws.Cells.Item(7, 12).Value = 500
$ws.Cells.Item(7, 13).Value = 475
$ws.Cells.Item(7, 16).Value = 475

$ws.Cells.Item(8, 4).Value = 44573
$ws.Cells.Item(8, 9).Value = "Tercera"
$ws.Cells.Item(8, 10).Value = 800
$ws.Cells.Item(8, 11).Value = 600
$ws.Cells.Item(8, 12).Value = 650
$ws.Cells.Item(8, 13).Value = 625
$ws.Cells.Item(8, 16).Value = 625

$ws.Cells.Item(9, 4).Value = 44210
$ws.Cells.Item(9, 9).Value = "Segunda"
$ws.Cells.Item(9, 10).Value = 900
$ws.Cells.Item(9, 11).Value = 600
$ws.Cells.Item(9, 12).Value = 700
$ws.Cells.Item(9, 13).Value = 650
$ws.Cells.Item(9, 16).Value = 650

$ws.Cells.Item(10, 4).Value = 44229
$ws.Cells.Item(10, 9).Value = "Segunda"
$ws.Cells.Item(10, 10).Value = 760
$ws.Cells.Item(10, 11).Value = 550
$ws.Cells.Item(10, 12).Value = 600
$ws.Cells.Item(10, 13).Value = 575
$ws.Cells.Item(10, 16).Value = 575

$ws.Cells.Item(11, 4).Value = 44278
$ws.Cells.Item(11, 9).Value = "Segunda"
$ws.Cells.Item(11, 10).Value = 700
$ws.Cells.Item(11, 11).Value = 600
$ws.Cells.Item(11, 12).Value = 700
$ws.Cells.Item(11, 13).Value = 650
$ws.Cells.Item(11, 16).Value = 650

$ws.Cells.Item(12, 4).Value = 44278
$ws.Cells.Item(12, 9).Value = "Tercera"
$ws.Cells.Item(12, 10).Value = 400
$ws.Cells.Item(12, 11).Value = 500
$ws.Cells.Item(12, 12).Value = 600
$ws.Cells.Item(12, 13).Value = 550
$ws.Cells.Item(12, 16).Value = 550

$ws.Cells.Item(13, 4).Value = 44658
$ws.Cells.Item(13, 9).Value = "Segunda"
$ws.Cells.Item(13, 10).Value = 1000
$ws.Cells.Item(13, 11).Value = 600
$ws.Cells.Item(13, 12).Value = 650
$ws.Cells.Item(13, 13).Value = 625
$ws.Cells.Item(13, 16).Value = 625

$ws.Cells.Item(14, 4).Value = 44253
$ws.Cells.Item(14, 9).Value = "Segunda"
$ws.Cells.Item(14, 10).Value = 1000
$ws.Cells.Item(14, 11).Value = 800
$ws.Cells.Item(14, 12).Value = 900
$ws.Cells.Item(14, 13).Value = 850
$ws.Cells.Item(14, 16).Value = 850

$ws.Cells.Item(15, 4).Value = 44253
$ws.Cells.Item(15, 9).Value = "Tercera"
$ws.Cells.Item(15, 10).Value = 800
$ws.Cells.Item(15, 11).Value = 600
$ws.Cells.Item(15, 12).Value = 700
$ws.Cells.Item(15, 13).Value = 650
$ws.Cells.Item(15, 16).Value = 650

$ws.Cells.Item(16, 4).Value = 44544
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 1000
$ws.Cells.Item(16, 11).Value = 600
$ws.Cells.Item(16, 12).Value = 650
$ws.Cells.Item(16, 13).Value = 625
$ws.Cells.Item(16, 16).Value = 625

$ws.Cells.Item(17, 4).Value = 44201
$ws.Cells.Item(17, 9).Value = "Segunda"
$ws.Cells.Item(17, 10).Value = 500
$ws.Cells.Item(17, 11).Value = 800
$ws.Cells.Item(17, 12).Value = 900
$ws.Cells.Item(17, 13).Value = 850
$ws.Cells.Item(17, 16).Value = 850

$ws.Cells.Item(18, 4).Value = 44174
$ws.Cells.Item(18, 9).Value = "Segunda"
$ws.Cells.Item(18, 10).Value = 800
$ws.Cells.Item(18, 11).Value = 450
$ws.Cells.Item(18, 12).Value = 500
$ws.Cells.Item(18, 13).Value = 475
$ws.Cells.Item(18, 16).Value = 475

$ws.Cells.Item(19, 4).Value = 44174
$ws.Cells.Item(19, 9).Value = "Tercera"
$ws.Cells.Item(19, 10).Value = 1200
$ws.Cells.Item(19, 11).Value = 250
$ws.Cells.Item(19, 12).Value = 350
$ws.Cells.Item(19, 13).Value = 300
$ws.Cells.Item(19, 16).Value = 300

$ws.Cells.Item(20, 4).Value = 44799
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 800
$ws.Cells.Item(20, 11).Value = 1000
$ws.Cells.Item(20, 12).Value = 1200
$ws.Cells.Item(20, 13).Value = 1100
$ws.Cells.Item(20, 16).Value = 1100

$ws.Cells.Item(21, 4).Value = 44267
$ws.Cells.Item(21, 9).Value = "Tercera"
$ws.Cells.Item(21, 10).Value = 400
$ws.Cells.Item(21, 11).Value = 500
$ws.Cells.Item(21, 12).Value = 600
$ws.Cells.Item(21, 13).Value = 550
$ws.Cells.Item(21, 16).Value = 550
